# Update cryptos list — GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 / Row 8 swap (USDC <-> XRP moved rank) ---
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"

# --- Row 19 / Row 20 swap (Uniswap <-> InternetComputer(DFINITY) moved rank) ---
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"

# --- Price (column D) values that are NOT ambiguous with numbers (kept as plain
#     text automatically because they contain multiple '.' separators or other
#     non-numeric characters) ---
$ws.Range("D2").Value = "43.043.96"
$ws.Range("D3").Value = "2.366.12"
$ws.Range("D15").Value = "2.733.26"
$ws.Range("D16").Value = "2.367.48"
$ws.Range("D18").Value = "43.018.69"
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("D43").Value = "1.928.20"
$ws.Range("D48").Value = "2.595.97"

# --- Price (column D) values that look like plain numbers: force text storage
#     (format as Text, assign, then restore the default "Normal" style so no
#     residual number formatting sticks to the cell) ---
$textPriceRefs = @(
    "D5","D6","D7","D8","D9","D10","D11","D12","D13","D14",
    "D17","D19","D20","D22","D23","D24","D25","D27","D28","D29",
    "D30","D32","D33","D34","D35","D36","D38","D39","D42","D44",
    "D46","D47","D49","D50","D51"
)
foreach ($ref in $textPriceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D5").Value  = "303.49"
$ws.Range("D6").Value  = "95.63"
$ws.Range("D7").Value  = "0.503"
$ws.Range("D8").Value  = "1.00"
$ws.Range("D9").Value  = "0.482"
$ws.Range("D10").Value = "34.35"
$ws.Range("D11").Value = "0.126"
$ws.Range("D12").Value = "0.0789"
$ws.Range("D13").Value = "18.30"
$ws.Range("D14").Value = "6.79"
$ws.Range("D17").Value = "0.803"
$ws.Range("D19").Value = "11.94"
$ws.Range("D20").Value = "6.30"
$ws.Range("D22").Value = "67.94"
$ws.Range("D23").Value = "235.34"
$ws.Range("D24").Value = "2.21"
$ws.Range("D25").Value = "2.44"
$ws.Range("D27").Value = "24.45"
$ws.Range("D28").Value = "2.37"
$ws.Range("D29").Value = "9.36"
$ws.Range("D30").Value = "32.07"
$ws.Range("D32").Value = "5.04"
$ws.Range("D33").Value = "0.111"
$ws.Range("D34").Value = "17.72"
$ws.Range("D35").Value = "0.0736"
$ws.Range("D36").Value = "129.02"
$ws.Range("D38").Value = "2.85"
$ws.Range("D39").Value = "4.32"
$ws.Range("D42").Value = "21.16"
$ws.Range("D44").Value = "0.0278"
$ws.Range("D46").Value = "2.76"
$ws.Range("D47").Value = "9.19"
$ws.Range("D49").Value = "1.52"
$ws.Range("D50").Value = "71.58"
$ws.Range("D51").Value = "51.61"

foreach ($ref in $textPriceRefs) {
    $ws.Range($ref).Style = "Normal"
}

# --- Volume(1h) (column E) values - padded percent strings, always stored as text ---
$ws.Range("E2").Value  = "  -0.12%  "
$ws.Range("E3").Value  = "  +1.10%  "
$ws.Range("E4").Value  = "  +0.07%  "
$ws.Range("E5").Value  = "  +0.22%  "
$ws.Range("E7").Value  = "  -0.14%  "
$ws.Range("E8").Value  = "  +0.01%  "
$ws.Range("E9").Value  = "  -2.59%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +4.39%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +10.39%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  +13.92%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  -8.47%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  -2.72%  "
